{"js": "// Apply the \"enchanted-eyes\" copy update: shorten the title, flesh out the\n// \"What we like\" / \"What we don't like\" bullet lists, and rewrite the final\n// bold title + italic summary line.\n\n// Replace every occurrence of `searchText` found via Body.search (used for\n// strings that are unique in the document, e.g. the title and the long\n// summary sentence).\nasync function replaceAllOccurrences(context, searchText, newText) {\n  const results = context.document.body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Replace the text of whichever paragraph's *entire* text exactly equals\n// `oldText`. This is used for the short bullet-list items, since a couple of\n// them (e.g. \"Egyptian-themed layout\") also occur as a substring elsewhere,\n// mid-sentence, and that unrelated text must not be touched.\nasync function replaceParagraphWithExactText(context, oldText, newText) {\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === oldText) {\n      paragraphs.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n  }\n  await context.sync();\n}\n\n// 1. Title (appears twice: the H1 heading and the bold line near the end).\nawait replaceAllOccurrences(\n  context,\n  \"Play Enchanted Eyes Online Slot Game for Free\",\n  \"Play Enchanted Eyes Slot Game Free\"\n);\n\n// 2. \"What we like\" bullets.\nawait replaceParagraphWithExactText(\n  context,\n  \"Stunning graphics\",\n  \"Stunning graphics and captivating Egyptian-themed layout\"\n);\nawait replaceParagraphWithExactText(\n  context,\n  \"Numerous bonuses and ways to win\",\n  \"Numerous bonuses and extra ways to win\"\n);\nawait replaceParagraphWithExactText(\n  context,\n  \"Egyptian-themed layout\",\n  \"Wild symbol (Princess) to help complete winning lines\"\n);\nawait replaceParagraphWithExactText(\n  context,\n  \"Wild and Scatter symbols\",\n  \"Chance to reactivate bonuses for more chances to win big\"\n);\n\n// 3. \"What we don't like\" bullets.\nawait replaceParagraphWithExactText(\n  context,\n  \"Low RTP at 94.31%\",\n  \"Limited availability of the Big Hit Bonanza games\"\n);\nawait replaceParagraphWithExactText(\n  context,\n  \"Bonus type received is random\",\n  \"Bonus symbol activation is random and based on bet level\"\n);\n\n// 4. Closing italic summary line.\nawait replaceAllOccurrences(\n  context,\n  \"Read our detailed review of Enchanted Eyes online slot game, play for free and discover numerous bonuses and ways to win in this stunning Egyptian-themed game.\",\n  \"Read our review of Enchanted Eyes, an online slot game with stunning graphics and numerous ways to win. Play for free now!\"\n);\n", "ps1": "# Apply the \"enchanted-eyes\" copy update: shorten the title, flesh out the\n# \"What we like\" / \"What we don't like\" bullet lists, and rewrite the final\n# bold title + italic summary line.\n\n$d = $word.ActiveDocument\n\n# Replace every occurrence of $oldText in the whole document (used for\n# strings that are unique/unambiguous, e.g. the title and the long summary\n# sentence that closes the document).\nfunction Replace-Everywhere($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n# Replace the text of whichever paragraph's *entire* text (sans the trailing\n# paragraph mark) exactly equals $oldText. This is used for the short\n# bullet-list items, since a couple of them (e.g. \"Egyptian-themed layout\")\n# also occur as a substring elsewhere, mid-sentence, or match another\n# heading case-insensitively (\"Egyptian-Themed Layout\"), and that unrelated\n# text must not be touched.\nfunction Replace-ExactParagraph($oldText, $newText) {\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t.CompareTo($oldText) -eq 0) {\n            $find = $p.Range.Find\n            $find.ClearFormatting()\n            $find.Replacement.ClearFormatting()\n            $find.Text = $oldText\n            $find.Replacement.Text = $newText\n            $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2) | Out-Null\n        }\n    }\n}\n\n# 1. Title (appears twice: the H1 heading and the bold line near the end).\nReplace-Everywhere \"Play Enchanted Eyes Online Slot Game for Free\" \"Play Enchanted Eyes Slot Game Free\"\n\n# 2. \"What we like\" bullets.\nReplace-ExactParagraph \"Stunning graphics\" \"Stunning graphics and captivating Egyptian-themed layout\"\nReplace-ExactParagraph \"Numerous bonuses and ways to win\" \"Numerous bonuses and extra ways to win\"\nReplace-ExactParagraph \"Egyptian-themed layout\" \"Wild symbol (Princess) to help complete winning lines\"\nReplace-ExactParagraph \"Wild and Scatter symbols\" \"Chance to reactivate bonuses for more chances to win big\"\n\n# 3. \"What we don't like\" bullets.\nReplace-ExactParagraph \"Low RTP at 94.31%\" \"Limited availability of the Big Hit Bonanza games\"\nReplace-ExactParagraph \"Bonus type received is random\" \"Bonus symbol activation is random and based on bet level\"\n\n# 4. Closing italic summary line.\nReplace-Everywhere \"Read our detailed review of Enchanted Eyes online slot game, play for free and discover numerous bonuses and ways to win in this stunning Egyptian-themed game.\" \"Read our review of Enchanted Eyes, an online slot game with stunning graphics and numerous ways to win. Play for free now!\"\n"}
